$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Text format on Price (D) and Volume(1h) (E) columns for the data rows
# so that numeric-looking strings (e.g. "322.48", "1.030") are kept as literal
# text instead of being auto-converted to numbers by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "27.722.74"
$ws.Range("E2").Value = "  +0.48%  "

$ws.Range("D3").Value = "1.853.71"
$ws.Range("E3").Value = "  +0.16%  "

$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").Value = "322.48"
$ws.Range("E5").Value = "  +0.29%  "

$ws.Range("D6").Value = "1.030"
$ws.Range("E6").Value = "  -0.17%  "

$ws.Range("D7").Value = "0.4395"
$ws.Range("E7").Value = "  +0.26%  "

$ws.Range("D8").Value = "0.3814"
$ws.Range("E8").Value = "  +1.64%  "

$ws.Range("D9").Value = "0.07422"
$ws.Range("E9").Value = "  +0.23%  "

$ws.Range("D10").Value = "0.8856"
$ws.Range("E10").Value = "  +1.07%  "

$ws.Range("D11").Value = "21.56"
$ws.Range("E11").Value = "  +0.39%  "

$ws.Range("D12").Value = "1.855.50"
$ws.Range("E12").Value = "  +0.06%  "

$ws.Range("D13").Value = "5.518"
$ws.Range("E13").Value = "  +0.04%  "

$ws.Range("D14").Value = "6.738"
$ws.Range("E14").Value = "  +0.65%  "

$ws.Range("D15").Value = "0.07200"
$ws.Range("E15").Value = "  -0.03%  "

$ws.Range("D16").Value = "85.93"
$ws.Range("E16").Value = "  +3.81%  "

$ws.Range("D17").Value = "1.037"
$ws.Range("E17").Value = "  +0.00%  "

$ws.Range("D18").Value = "0.000009088"
$ws.Range("E18").Value = "  +0.61%  "

$ws.Range("D19").Value = "1.031"
$ws.Range("E19").Value = "  +0.00%  "

$ws.Range("D20").Value = "15.54"
$ws.Range("E20").Value = "  +0.57%  "

$ws.Range("D21").Value = "27.726.31"
$ws.Range("E21").Value = "  +0.47%  "

$ws.Range("D22").Value = "5.292"
$ws.Range("E22").Value = "  +0.48%  "

$ws.Range("D23").Value = "11.26"
$ws.Range("E23").Value = "  +0.15%  "

$ws.Range("D24").Value = "2.084.44"
$ws.Range("E24").Value = "  +0.50%  "

$ws.Range("D25").Value = "2.069"
$ws.Range("E25").Value = "  +5.98%  "

$ws.Range("D26").Value = "158.85"
$ws.Range("E26").Value = "  +0.63%  "

$ws.Range("D27").Value = "18.73"
$ws.Range("E27").Value = "  -0.01%  "

$ws.Range("D28").Value = "5.344"
$ws.Range("E28").Value = "  +0.76%  "

$ws.Range("D29").Value = "1.986"
$ws.Range("E29").Value = "  +2.31%  "

$ws.Range("D30").Value = "118.61"
$ws.Range("E30").Value = "  +1.98%  "

$ws.Range("D31").Value = "0.09105"
$ws.Range("E31").Value = "  +0.50%  "

$ws.Range("D32").Value = "0.7719"
$ws.Range("E32").Value = "  +0.55%  "

$ws.Range("D33").Value = "1.209"
$ws.Range("E33").Value = "  +0.07%  "

$ws.Range("D34").Value = "3.039"
$ws.Range("E34").Value = "  +5.04%  "

$ws.Range("D35").Value = "4.592"
$ws.Range("E35").Value = "  +1.38%  "

$ws.Range("D36").Value = "1.033"
$ws.Range("E36").Value = "  -0.24%  "

$ws.Range("D37").Value = "1.151"
$ws.Range("E37").Value = "  -0.30%  "

$ws.Range("D38").Value = "0.01982"
$ws.Range("E38").Value = "  +0.19%  "

$ws.Range("D39").Value = "0.05308"
$ws.Range("E39").Value = "  +0.37%  "

$ws.Range("D40").Value = "2.856"
$ws.Range("E40").Value = "  +1.22%  "

$ws.Range("D41").Value = "0.5199"
$ws.Range("E41").Value = "  +0.51%  "

$ws.Range("D42").Value = "6.951"
$ws.Range("E42").Value = "  +3.15%  "

$ws.Range("D43").Value = "0.1675"
$ws.Range("E43").Value = "  +0.12%  "

$ws.Range("D44").Value = "8.768"
$ws.Range("E44").Value = "  +2.12%  "

$ws.Range("D45").Value = "10.78"
$ws.Range("E45").Value = "  +1.78%  "

$ws.Range("D46").Value = "110.14"
$ws.Range("E46").Value = "  +1.12%  "

$ws.Range("D47").Value = "1.033"
$ws.Range("E47").Value = "  -0.14%  "

$ws.Range("D48").Value = "0.06551"
$ws.Range("E48").Value = "  +2.49%  "

$ws.Range("D49").Value = "1.710"
$ws.Range("E49").Value = "  -0.29%  "

$ws.Range("D50").Value = "0.4725"
$ws.Range("E50").Value = "  +1.42%  "

$ws.Range("D51").Value = "1.884"
$ws.Range("E51").Value = "  -0.57%  "
